# Auto-generated edit script for BRVM Recommandations workbook update
# Updates sector rows (2-11), reorders/updates stock rows (15-38) with 3 new rows (36-38)
# on sheet 'Recommandations', and updates the 'Top_YTD' sheet (rows 2-11).

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

function Set-RecoRow($row, $a, $b, $c, $d, $e, $f, $g) {
  $wsReco.Cells.Item($row, 1).Value = $a
  $wsReco.Cells.Item($row, 2).Value = $b
  $wsReco.Cells.Item($row, 3).Value = $c
  $wsReco.Cells.Item($row, 4).Value = $d
  $wsReco.Cells.Item($row, 5).Value = $e
  $wsReco.Cells.Item($row, 6).Value = $f
  $wsReco.Cells.Item($row, 7).Value = $g
}

function Set-YtdRow($row, $a, $b) {
  $wsYtd.Cells.Item($row, 1).Value = $a
  $wsYtd.Cells.Item($row, 2).Value = $b
}

# --- Recommandations sheet ---
Set-RecoRow 2 'BRVM - CONSOMMATION DISCRETIONNAIRE' 0 4 681.2 167.98 '🟡 Observer' '➖ Neutre'
Set-RecoRow 3 'BRVM - SERVICES FINANCIERS' 0 4 598.29 150.21 '🟡 Observer' '➖ Neutre'
Set-RecoRow 4 'BRVM - INDUSTRIELS' 0 4 583.34 145.87 '🟡 Observer' '➖ Neutre'
Set-RecoRow 5 'BRVM-PRESTIGE' 0 4 580.0599999999999 145.81 '🟡 Observer' '➖ Neutre'
Set-RecoRow 6 'BRVM - SERVICES PUBLICS' 0 4 461.18 117.36 '🟡 Observer' '➖ Neutre'
Set-RecoRow 7 'BRVM - CONSOMMATION DE BASE     (**)' 0 2 460.08 231.36 '🟡 Observer' '➖ Neutre'
Set-RecoRow 8 'BRVM - ENERGIE' 0 4 458.49 116.07 '🟡 Observer' '➖ Neutre'
Set-RecoRow 9 'BRVM-PRINCIPAL     (**)' 0 2 456 228.7 '🟡 Observer' '➖ Neutre'
Set-RecoRow 10 'BRVM - TELECOMMUNICATIONS' 0 4 377.03 94.59999999999999 '🟡 Observer' '➖ Neutre'
Set-RecoRow 11 'BRVM – COMPOSITE TOTAL RETURN     (**)' 0 2 271.14 136 '🟡 Observer' '➖ Neutre'
Set-RecoRow 15 'BERNABE CI (BNBC)' 3 0 22.05 7.4 '🟢 Achat' '✅ Renforcer'
Set-RecoRow 16 'SICABLE CI (CABC)' 2 0 11.95 7.34 '🟡 Observer' '➖ Neutre'
Set-RecoRow 17 'UNILEVER CI (UNLC)' 2 1 8.07 7.49 '🟡 Observer' '👀 À surveiller'
Set-RecoRow 18 'EVIOSYS PACKAGING SIEM CI (SEMC)' 2 1 7.06 -6.81 '🟢 Achat' '👀 À surveiller'
Set-RecoRow 19 'SAFCA CI (SAFC)' 1 0 5.74 5.74 '🟡 Observer' '➖ Neutre'
Set-RecoRow 20 'NSIA BANQUE COTE D''IVOIRE (NSBC)' 1 0 5.04 5.04 '🟡 Observer' '➖ Neutre'
Set-RecoRow 21 'ONATEL BF (ONTBF)' 1 0 4.38 4.38 '🟡 Observer' '➖ Neutre'
Set-RecoRow 22 'SITAB CI (STBC)' 1 0 3.68 3.68 '🟡 Observer' '➖ Neutre'
Set-RecoRow 23 'NESTLE CI (NTLC)' 1 0 3.64 3.64 '🟡 Observer' '➖ Neutre'
Set-RecoRow 24 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)' 1 0 3.35 3.35 '🟡 Observer' '➖ Neutre'
Set-RecoRow 25 'CORIS BANK INTERNATIONAL (CBIBF)' 1 1 3.22 7.11 '🟡 Observer' '👀 À surveiller'
Set-RecoRow 26 'SMB CI (SMBC)' 1 0 2.99 2.99 '🟡 Observer' '➖ Neutre'
Set-RecoRow 27 'TRACTAFRIC MOTORS CI (PRSC)' 1 1 2.44 -2.81 '🟡 Observer' '👀 À surveiller'
Set-RecoRow 28 'SOLIBRA CI (SLBC)' 1 1 0.86 -5.22 '🟡 Observer' '👀 À surveiller'
Set-RecoRow 29 'BICI CI (BICC)' 0 1 -0.49 -0.49 '🟡 Observer' '➖ Neutre'
Set-RecoRow 30 'TOTALENERGIES MARKETING CI (TTLC)' 0 1 -1.88 -1.88 '🟡 Observer' '➖ Neutre'
Set-RecoRow 31 'TOTALENERGIES MARKETING SN (TTLS)' 0 1 -1.92 -1.92 '🟡 Observer' '➖ Neutre'
Set-RecoRow 32 'SOGB CI (SOGC)' 0 1 -2.04 -2.04 '🟡 Observer' '➖ Neutre'
Set-RecoRow 33 'ORAGROUP TOGO (ORGT)' 0 1 -2.08 -2.08 '🟡 Observer' '➖ Neutre'
Set-RecoRow 34 'SETAO CI (STAC)' 0 1 -4 -4 '🟡 Observer' '➖ Neutre'
Set-RecoRow 35 'ECOBANK TRANS. INCORP. TG (ETIT)' 0 2 -4.15 -4.35 '🟡 Observer' '👀 À surveiller'
Set-RecoRow 36 'FILTISAC CI (FTSC)' 0 2 -4.37 -3.67 '🟡 Observer' '➖ Neutre'
Set-RecoRow 37 'CFAO MOTORS CI (CFAC)' 0 3 -4.9 -1.43 '🟡 Observer' '➖ Neutre'
Set-RecoRow 38 'AFRICA GLOBAL LOGISTICS CI (SDSC)' 0 2 -5.7 -2.89 '🟡 Observer' '➖ Neutre'

# --- Top_YTD sheet ---
Set-YtdRow 2 'BRVM - CONSOMMATION DISCRETIONNAIRE' 5237.65
Set-YtdRow 3 'BRVM - SERVICES FINANCIERS' 3779.51
Set-YtdRow 4 'BRVM - INDUSTRIELS' 3552.36
Set-YtdRow 5 'BRVM-PRESTIGE' 3503.83
Set-YtdRow 6 'BRVM - SERVICES PUBLICS' 2048.28
Set-YtdRow 7 'BRVM - ENERGIE' 2021.69
Set-YtdRow 8 'BRVM - TELECOMMUNICATIONS' 1324
Set-YtdRow 9 'BRVM - CONSOMMATION DE BASE     (**)' 989.25
Set-YtdRow 10 'BRVM-PRINCIPAL     (**)' 975.84
Set-YtdRow 11 'BRVM – COMPOSITE TOTAL RETURN     (**)' 454.93
